$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
}

# Row 2 - Bitcoin
Set-TextValue "D2" "63.093.87"
Set-TextValue "E2" "  +2.65%  "

# Row 3 - Ethereum
Set-TextValue "D3" "3.483.28"
Set-TextValue "E3" "  +2.62%  "

# Row 4 - TetherUSD
Set-TextValue "E4" "  +0.05%  "

# Row 5 - BNB
Set-TextValue "D5" "584.35"
Set-TextValue "E5" "  +1.36%  "

# Row 6 - Solana
Set-TextValue "D6" "147.93"
Set-TextValue "E6" "  +5.08%  "

# Row 8 - XRP
Set-TextValue "D8" "0.479"
Set-TextValue "E8" "  +1.03%  "

# Row 9 - Toncoin
Set-TextValue "E9" "  -0.07%  "

# Row 10 - Dogecoin
Set-TextValue "E10" "  +2.28%  "

# Row 11 - Cardano
Set-TextValue "E11" "  +2.73%  "

# Row 12 - WrappedliquidstakedEther2.0
Set-TextValue "D12" "4.086.06"
Set-TextValue "E12" "  +2.77%  "

# Row 13 - Avalanche
Set-TextValue "D13" "29.85"
Set-TextValue "E13" "  +5.81%  "

# Row 14 - TRON
Set-TextValue "E14" "  -0.13%  "

# Row 15 - WrappedEther
Set-TextValue "D15" "3.495.28"
Set-TextValue "E15" "  +2.87%  "

# Row 16 - ShibaInu
Set-TextValue "E16" "  +1.24%  "

# Row 17 - WrappedBTC
Set-TextValue "D17" "63.147.60"
Set-TextValue "E17" "  +2.70%  "

# Row 18 - Polkadot
Set-TextValue "D18" "6.32"
Set-TextValue "E18" "  +2.99%  "

# Row 19 - Chainlink
Set-TextValue "D19" "14.36"
Set-TextValue "E19" "  +4.95%  "

# Row 20 - Uniswap
Set-TextValue "E20" "  +4.37%  "

# Row 21 - BitcoinCash
Set-TextValue "D21" "388.89"
Set-TextValue "E21" "  -0.60%  "

# Row 22 - Polygon
Set-TextValue "D22" "0.564"
Set-TextValue "E22" "  +1.67%  "

# Row 23 - Litecoin
Set-TextValue "D23" "75.06"
Set-TextValue "E23" "  -0.45%  "

# Row 24 - Dai
Set-TextValue "E24" "  -0.06%  "

# Row 25 - swapped from PEPE to WrappedeETH
Set-TextValue "B25" "WrappedeETH"
Set-TextValue "C25" "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
Set-TextValue "D25" "3.631.83"
Set-TextValue "E25" "  +2.89%  "

# Row 26 - swapped from WrappedeETH to PEPE
Set-TextValue "B26" "PEPE"
Set-TextValue "C26" "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-TextValue "D26" "0.0000117"
Set-TextValue "E26" "  +3.34%  "

# Row 27 - Kaspa
Set-TextValue "D27" "0.181"
Set-TextValue "E27" "  -6.25%  "

# Row 28 - RenderToken
Set-TextValue "D28" "7.70"
Set-TextValue "E28" "  +5.57%  "

# Row 29 - Binance-PegBSC-USD
Set-TextValue "D29" "1.00"
Set-TextValue "E29" "  +0.03%  "

# Row 30 - InternetComputer(DFINITY)
Set-TextValue "D30" "8.27"
Set-TextValue "E30" "  +2.96%  "

# Row 31 - swapped from Fetch.AI to PancakeSwap
Set-TextValue "B31" "PancakeSwap"
Set-TextValue "C31" "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue "D31" "2.14"
Set-TextValue "E31" "  -0.03%  "

# Row 32 - swapped from PancakeSwap to Fetch.AI
Set-TextValue "B32" "Fetch.AI"
Set-TextValue "C32" "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue "D32" "1.43"
Set-TextValue "E32" "  +4.10%  "

# Row 33 - USDe
Set-TextValue "E33" "  -0.02%  "

# Row 34 - EthereumClassic
Set-TextValue "D34" "23.80"
Set-TextValue "E34" "  +1.69%  "

# Row 35 - NEARProtocol
Set-TextValue "D35" "5.33"
Set-TextValue "E35" "  +5.62%  "

# Row 36 - Aptos
Set-TextValue "D36" "7.12"
Set-TextValue "E36" "  +2.61%  "

# Row 37 - EnergySwap
Set-TextValue "D37" "31.93"
Set-TextValue "E37" "  +22.34%  "

# Row 38 - ImmutableX
Set-TextValue "E38" "  +6.98%  "

# Row 39 - Monero
Set-TextValue "D39" "171.17"
Set-TextValue "E39" "  +2.15%  "

# Row 40 - RenzoRestakedETH
Set-TextValue "D40" "3.522.80"
Set-TextValue "E40" "  +2.78%  "

# Row 41 - Hedera
Set-TextValue "E41" "  +0.55%  "

# Row 42 - Mantle
Set-TextValue "D42" "0.808"
Set-TextValue "E42" "  +3.74%  "

# Row 43 - OKB
Set-TextValue "D43" "42.37"
Set-TextValue "E43" "  -0.14%  "

# Row 44 - Filecoin
Set-TextValue "D44" "4.48"
Set-TextValue "E44" "  +1.23%  "

# Row 45 - Stacks
Set-TextValue "D45" "1.72"
Set-TextValue "E45" "  +3.57%  "

# Row 46 - ONDO
Set-TextValue "E46" "  +6.05%  "

# Row 47 - Maker
Set-TextValue "D47" "2.629.59"
Set-TextValue "E47" "  +6.77%  "

# Row 48 - dogwifhat
Set-TextValue "D48" "2.28"
Set-TextValue "E48" "  +10.13%  "

# Row 49 - InjectiveProtocol
Set-TextValue "D49" "23.42"
Set-TextValue "E49" "  +1.74%  "

# Row 50 - Cosmos
Set-TextValue "E50" "  +1.19%  "

# Row 51 - VeChain
Set-TextValue "E51" "  +2.69%  "
